$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Equity Portfolio
$ws.Range("A2").Value = "Equity Portfolio"
$ws.Range("B2").Value = 127.137217321074
$ws.Range("C2").Value = 14.76687513140331
$ws.Range("D2").Value = 24.05919947431317
$ws.Range("E2").Value = 0.5306442196895986
$ws.Range("F2").Value = -34.39828036199673

# Row 3 - Covered Call Strategy
$ws.Range("A3").Value = "Covered Call Strategy"
$ws.Range("B3").Value = 104.4209511983805
$ws.Range("C3").Value = 12.75440364326808
$ws.Range("D3").Value = 17.92061497814814
$ws.Range("E3").Value = 0.6001135371962217
$ws.Range("F3").Value = -30.60266810485754

# Row 4 - Combined Portfolio
$ws.Range("A4").Value = "Combined Portfolio"
$ws.Range("B4").Value = 123.5200755309552
$ws.Range("C4").Value = 14.45798093377535
$ws.Range("D4").Value = 21.18622196611397
$ws.Range("E4").Value = 0.5880227703505186
$ws.Range("F4").Value = -33.1950418168004

# Row 5 - SPY Buy & Hold (values only, label unchanged)
$ws.Range("B5").Value = 125.1812014728365
$ws.Range("C5").Value = 14.6003491757158
$ws.Range("D5").Value = 19.70244352451648
$ws.Range("E5").Value = 0.6395323077585131
$ws.Range("F5").Value = -33.71726063766723

$wb.Save()
